# Done chức năng bán hàng + ghi hoá đơn
# Renumber "Mã SP" (product codes) to a contiguous sequence 111-132,
# fix the "Trà sữa trà đen" quantity, and drop the discontinued
# "Trà đào cam sả" drink row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contiguous product codes for rows 2..23 (Mã SP column A)
$newCodes = @("111","112","113","114","115","116","117","118","119","120","121","122","123","124","125","126","127","128","129","130","131","132")

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newCodes[$i]
}

# "Trà sữa trà đen" (row 23) Số lượng corrected from 9 to 25
$ws.Range("E23").Value = "25"

# Remove the discontinued "Trà đào cam sả" row (old row 24) entirely
$ws.Rows.Item(24).Delete()
